$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5, shifting existing rows 5-17 down to 6-18
$ws.Rows.Item(5).Insert()

# Fill in the new row 5 with data (new weekly record dated 2021-11-25 / serial 44525)
$ws.Cells.Item(5, 1).Value = 7
$ws.Cells.Item(5, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(5, 3).Value = "Ñuble"
$ws.Cells.Item(5, 4).Value = 44525
$ws.Cells.Item(5, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5, 5).Value = 16
$ws.Cells.Item(5, 6).Value = 100112040
$ws.Cells.Item(5, 7).Value = "Cilantro"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 60
$ws.Cells.Item(5, 11).Value = 2000
$ws.Cells.Item(5, 12).Value = 2000
$ws.Cells.Item(5, 13).Value = 2000
$ws.Cells.Item(5, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(5, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(5, 16).Value = 2000
$ws.Cells.Item(5, 17).Value = 1
$ws.Cells.Item(5, 18).Value = "Hortaliza"
